$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.637
$ws.Range("D5").Value = 0.723
$ws.Range("E5").Value = 0.744
$ws.Range("F5").Value = 0.77
$ws.Range("G5").Value = 0.651
$ws.Range("H5").Value = 0.662

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.637
$ws.Range("E7").Value = 0.744
$ws.Range("F7").Value = 0.77
$ws.Range("H7").Value = 0.662

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.658
$ws.Range("D8").Value = 0.767
$ws.Range("E8").Value = 0.788
$ws.Range("F8").Value = 0.803
$ws.Range("G8").Value = 0.719
$ws.Range("H8").Value = 0.734

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.579
$ws.Range("C9").Value = 0.679
$ws.Range("D9").Value = 0.774
$ws.Range("E9").Value = 0.789
$ws.Range("F9").Value = 0.8
$ws.Range("G9").Value = 0.709
$ws.Range("H9").Value = 0.719

$wb.Save()
